$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right
#    after the first "Play Burning Slots Cash Mesh Free - Review" heading.
# ----------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete() | Out-Null
        break
    }
}

# ----------------------------------------------------------------------
# 2) Insert a new paragraph - "Play Burning Slots Cash Mesh Free - Review"
#    in bold - right before the final (now FAQ-text) paragraph.
# ----------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$breakPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$breakPoint.InsertParagraphBefore() | Out-Null

$newParaIndex = $count
$newPara = $d.Paragraphs.Item($newParaIndex)
$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$newParaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Burning Slots Cash Mesh Free - Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newParaRange.InsertXML($newParaXml) | Out-Null

# ----------------------------------------------------------------------
# 3) Replace the FAQ text in the last paragraph with the meta-description
#    wording (keeping its existing italic formatting).
# ----------------------------------------------------------------------
$oldText = "Can I play " + [char]34 + "Burning Slots Cash Mesh" + [char]34 + " without registration? Yes, you can try out the demo version without registration. What is the maximum payout for this game? The maximum payout is 50,000x your bet. What is the minimum and maximum bet amount? The minimum bet is " + [char]0x20AC + "0.10, and the maximum bet amount is " + [char]0x20AC + "50. How many paylines does " + [char]34 + "Burning Slots Cash Mesh" + [char]34 + " have? It has five paylines."
$newText = "Read our review of Burning Slots Cash Mesh and play for free. Features Cash Mesh and Respins for high payouts."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newText, 2) | Out-Null
